$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new value (applies to both column B and D, stored as text
# to match the original "0.00"-style formatted numeric strings).
$updates = @{
    2 = "15.00"
    3 = "1.00"
    4 = "8.00"
    5 = "4.00"
    6 = "4.00"
    7 = "32.00"
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]

    $rangeB = $ws.Range("B$row")
    $rangeB.NumberFormat = "@"
    $rangeB.Value = $value

    $rangeD = $ws.Range("D$row")
    $rangeD.NumberFormat = "@"
    $rangeD.Value = $value
}
